$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("01142019")
$ws.Rows.Item(15).Select()
$ws.Rows.Item(15).Delete()
